$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "-"
$ws.Range("E2").Value = "MCT-3A-Robótica"

# Row 3
$ws.Range("B3").Value = "-"
$ws.Range("E3").Value = "MCT-3A-Robótica"
$ws.Range("F3").Value = "MCT-3A-Robótica"

# Row 4
$ws.Range("D4").Value = "-"
$ws.Range("F4").Value = "MCT-3A-Robótica"

# Row 6
$ws.Range("D6").Value = "-"

# Row 10
$ws.Range("E10").Value = "MEC-2A-Elet. Digit. Básica"

# Row 11
$ws.Range("E11").Value = "MEC-2A-Elet. Digit. Básica"

# Row 12
$ws.Range("B12").Value = "MEC-2A-Elet. Digit. Básica"
$ws.Range("E12").Value = "MEC-2A-Elet. Digit. Básica"
$ws.Range("F12").Value = "-"

# Row 14
$ws.Range("B14").Value = "MEC-2A-Elet. Digit. Básica"
$ws.Range("E14").Value = "MEC-2A-Elet. Digit. Básica"
$ws.Range("F14").Value = "-"

# Row 15
$ws.Range("E15").Value = "MEC-2A-Elet. Digit. Básica"
$ws.Range("F15").Value = "-"

# Row 16
$ws.Range("E16").Value = "MEC-2A-Elet. Digit. Básica"
$ws.Range("F16").Value = "-"
